# Apply updated crypto price / volume(1h) data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.494.70"

$ws.Cells.Item(3, 4).Value = "3.257.09"
$ws.Cells.Item(3, 5).Value = "  -0.58%  "

$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 5).Value = "  -0.15%  "

$ws.Cells.Item(5, 4).Value = "'574.86"
$ws.Cells.Item(5, 5).Value = "  -0.43%  "

$ws.Cells.Item(6, 4).Value = "'177.88"
$ws.Cells.Item(6, 5).Value = "  -3.19%  "

$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "  +0.00%  "

$ws.Cells.Item(8, 4).Value = "'0.581"
$ws.Cells.Item(8, 5).Value = "  +2.34%  "

$ws.Cells.Item(9, 4).Value = "3.253.15"
$ws.Cells.Item(9, 5).Value = "  -0.62%  "

$ws.Cells.Item(10, 4).Value = "'0.174"
$ws.Cells.Item(10, 5).Value = "  -0.12%  "

$ws.Cells.Item(11, 4).Value = "'0.570"
$ws.Cells.Item(11, 5).Value = "  +0.14%  "

$ws.Cells.Item(12, 4).Value = "'45.27"
$ws.Cells.Item(12, 5).Value = "  -1.73%  "

$ws.Cells.Item(13, 5).Value = "  +2.19%  "

$ws.Cells.Item(14, 4).Value = "'677.81"
$ws.Cells.Item(14, 5).Value = "  +10.64%  "

$ws.Cells.Item(15, 4).Value = "3.791.54"
$ws.Cells.Item(15, 5).Value = "  -0.38%  "

$ws.Cells.Item(16, 4).Value = "'8.30"
$ws.Cells.Item(16, 5).Value = "  -1.10%  "

$ws.Cells.Item(17, 4).Value = "67.533.66"
$ws.Cells.Item(17, 5).Value = "  +3.02%  "

$ws.Cells.Item(18, 5).Value = "  +1.34%  "

$ws.Cells.Item(19, 4).Value = "3.271.55"
$ws.Cells.Item(19, 5).Value = "  -0.27%  "

$ws.Cells.Item(20, 4).Value = "'17.27"
$ws.Cells.Item(20, 5).Value = "  -2.75%  "

$ws.Cells.Item(21, 4).Value = "'10.69"
$ws.Cells.Item(21, 5).Value = "  -1.74%  "

$ws.Cells.Item(22, 4).Value = "'0.885"
$ws.Cells.Item(22, 5).Value = "  -0.13%  "

$ws.Cells.Item(23, 4).Value = "'17.01"
$ws.Cells.Item(23, 5).Value = "  -5.36%  "

$ws.Cells.Item(24, 4).Value = "'5.09"
$ws.Cells.Item(24, 5).Value = "  +2.91%  "

$ws.Cells.Item(25, 4).Value = "'97.98"
$ws.Cells.Item(25, 5).Value = "  -2.37%  "

$ws.Cells.Item(26, 5).Value = "  -2.37%  "

$ws.Cells.Item(27, 5).Value = "  +0.88%  "

$ws.Cells.Item(28, 4).Value = "'9.31"
$ws.Cells.Item(28, 5).Value = "  -1.23%  "

$ws.Cells.Item(29, 4).Value = "'32.46"
$ws.Cells.Item(29, 5).Value = "  +5.63%  "

$ws.Cells.Item(30, 4).Value = "'8.37"
$ws.Cells.Item(30, 5).Value = "  -0.37%  "

$ws.Cells.Item(31, 4).Value = "'6.62"
$ws.Cells.Item(31, 5).Value = "  +3.36%  "

$ws.Cells.Item(32, 4).Value = "'581.57"
$ws.Cells.Item(32, 5).Value = "  +6.27%  "

$ws.Cells.Item(33, 4).Value = "3.858.48"
$ws.Cells.Item(33, 5).Value = "  +2.34%  "

$ws.Cells.Item(34, 4).Value = "'10.78"
$ws.Cells.Item(34, 5).Value = "  -0.27%  "

$ws.Cells.Item(35, 5).Value = "  +0.18%  "

$ws.Cells.Item(36, 5).Value = "  -0.07%  "

$ws.Cells.Item(37, 4).Value = "'3.34"
$ws.Cells.Item(37, 5).Value = "  -9.49%  "

$ws.Cells.Item(38, 4).Value = "'55.13"
$ws.Cells.Item(38, 5).Value = "  -1.55%  "

$ws.Cells.Item(39, 5).Value = "  +1.30%  "

$ws.Cells.Item(40, 4).Value = "'3.21"
$ws.Cells.Item(40, 5).Value = "  +2.94%  "

$ws.Cells.Item(41, 2).Value = "ApeXProtocol"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(41, 4).Value = "'3.41"
$ws.Cells.Item(41, 5).Value = "  +0.89%  "

$ws.Cells.Item(42, 2).Value = "Fetch.AI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(42, 4).Value = "'2.61"
$ws.Cells.Item(42, 5).Value = "  +1.41%  "

$ws.Cells.Item(43, 4).Value = "'31.97"
$ws.Cells.Item(43, 5).Value = "  -1.08%  "

$ws.Cells.Item(44, 4).Value = "0.0₃0666"
$ws.Cells.Item(44, 5).Value = "  -0.97%  "

$ws.Cells.Item(45, 4).Value = "'0.328"

$ws.Cells.Item(46, 4).Value = "'0.0409"
$ws.Cells.Item(46, 5).Value = "  +1.27%  "

$ws.Cells.Item(47, 5).Value = "  +1.42%  "

$ws.Cells.Item(48, 5).Value = "  +0.35%  "

$ws.Cells.Item(49, 5).Value = "  +9.24%  "

$ws.Cells.Item(50, 5).Value = "  -0.22%  "

$ws.Cells.Item(51, 4).Value = "'129.67"
$ws.Cells.Item(51, 5).Value = "  +1.03%  "
